$d = $word.ActiveDocument

# --- Edit 1: paragraph "Ha : At least two means are different" -------------
# Originally the trailing ":" and "At least two means are different" text
# live in two separate <m:r> math runs (after the H_a <m:sSub>). Merge them
# into a single run, keeping the <m:sSub> for "Ha" untouched.
$p1 = $d.Paragraphs.Item(57)
$om1 = $p1.Range.OMaths.Item(1)
$xml1 = '<m:oMathPara><m:oMath>' + `
  '<m:sSub>' + `
    '<m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr>' + `
    '<m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>H</m:t></m:r></m:e>' + `
    '<m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>a</m:t></m:r></m:sub>' + `
  '</m:sSub>' + `
  '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>:At least two means are different</m:t></m:r>' + `
  '</m:oMath></m:oMathPara>'
$om1.Range.InsertXML($xml1)

# --- Edit 2: paragraph "p-value=0.0024" -------------------------------------
# Originally "p", "-value=", and "0.0024" live in three separate <m:r> math
# runs. Merge them into a single run.
$p2 = $d.Paragraphs.Item(78)
$om2 = $p2.Range.OMaths.Item(1)
$xml2 = '<m:oMathPara><m:oMath>' + `
  '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>p-value=0.0024</m:t></m:r>' + `
  '</m:oMath></m:oMathPara>'
$om2.Range.InsertXML($xml2)

Write-Output "done"
